# Applies the edits described by the commit:
#   1. Notes master "datetimeFigureOut" footer field: 15-12-2023 -> 21-12-2023
#   2. Title-slide layout headline: merge the "Digital " + "Adventure " runs
#      into a single run "Digital Adventure "
#   3. "Thank You" layout textbox: add baseline="0" and recolor FE04AC -> C00000

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Notes Master date placeholder field text
# ---------------------------------------------------------------------------
$notesMaster = $p.NotesMaster
$dateAndTime = $notesMaster.HeadersFooters.DateAndTime
$dateAndTime.Text = "21-12-2023"

# ---------------------------------------------------------------------------
# 2) Slide Layout 1 ("Digital Adventure Ride to the Future" title) - merge runs
# ---------------------------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$layouts = $master.CustomLayouts

$titleLayout = $layouts.Item(1)
$titleShape = $titleLayout.Shapes.Item(3)
$titleShape.TextFrame.TextRange.Text = "Digital Adventure `rRide to the Future"

# ---------------------------------------------------------------------------
# 3) Slide Layout 4 ("Thank You" layout) - baseline + color
# ---------------------------------------------------------------------------
$thankYouLayout = $layouts.Item(4)
$thankYouShape = $thankYouLayout.Shapes.Item(2)
$thankYouFont = $thankYouShape.TextFrame.TextRange.Font
$thankYouFont.BaselineOffset = 0
$thankYouFont.Color.RGB = 192
